$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.124.75'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.97%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.750.91'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.47%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.48%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5279'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.74%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2809'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.48%  '

$ws.Range("E9").Value = '  +1.86%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.746.68'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.18%  '

$ws.Range("E11").Value = '  +2.85%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.51'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.93%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6472'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.19%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.633'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.26%  '

$ws.Range("E15").Value = '  +2.81%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9999'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.09%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9997'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.04%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.017.23'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.77'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.84%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006745'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.57%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.969.10'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.342'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.17%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.755'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.95%  '

$ws.Range("E24").Value = '  +2.63%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '138.95'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.26%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.518'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.55%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.35'
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.815'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.37%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '104.94'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.90%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08294'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.15%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.817'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.68%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.660'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.44%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04615'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.46%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.645'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.014'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.92%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6366'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.05%  '

$ws.Range("E37").Value = '  +1.34%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01612'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.53%  '

$ws.Range("E39").Value = '  +4.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9994'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.19%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '101.73'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.20%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.3958'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.52%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7465'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.53%  '

$ws.Range("E44").Value = '  +2.39%  '

$ws.Range("E45").Value = '  +4.53%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.397'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.91%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05346'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.83%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.46'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.93%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.99'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.55%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3489'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.33%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.601'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.88%  '
